$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Row 11 (Marking): Right count 5 -> 4, Wrong count -1 -> -2
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Row 12 (Total): Right total 115 -> 92, Wrong total -3 -> -6, Max string "115 / 140" -> "86 / 112"
$ws.Range("B12").Value = 92
$ws.Range("C12").Value = -6
$ws.Range("E12").Value = "86 / 112"
